# Auto-generated Excel COM-interop edit script.
# Source: diff updating '江西-漫展信息.xlsx' (commit 456a3b4).
#
# 1) Sheet '展览'    : refresh column F (想去人数) values.
# 2) Sheet '全部类型' : drop the 4 duplicate trailing rows
#    (dimension A1:I40 -> A1:I36) and refresh the surviving
#    36 rows (incl. column F) to match '展览'.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$fExpo = New-Object 'object[,]' 35,1
$fExpo[0,0] = 36
$fExpo[1,0] = 197
$fExpo[2,0] = 50
$fExpo[3,0] = 250
$fExpo[4,0] = 40
$fExpo[5,0] = 138
$fExpo[6,0] = 253
$fExpo[7,0] = 12
$fExpo[8,0] = 249
$fExpo[9,0] = 14
$fExpo[10,0] = 37
$fExpo[11,0] = 25
$fExpo[12,0] = 86
$fExpo[13,0] = 433
$fExpo[14,0] = 44
$fExpo[15,0] = 476
$fExpo[16,0] = 402
$fExpo[17,0] = 137
$fExpo[18,0] = 64
$fExpo[19,0] = 32
$fExpo[20,0] = 38
$fExpo[21,0] = 1099
$fExpo[22,0] = 2832
$fExpo[23,0] = 22
$fExpo[24,0] = 55
$fExpo[25,0] = 536
$fExpo[26,0] = 46
$fExpo[27,0] = 1608
$fExpo[28,0] = 566
$fExpo[29,0] = 452
$fExpo[30,0] = 261
$fExpo[31,0] = 390
$fExpo[32,0] = 450
$fExpo[33,0] = 600
$fExpo[34,0] = 422
$wsExpo.Range("F2:F36").Value = $fExpo

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows("37:40").Delete()

$wsAll.Range("B1:B36").NumberFormat = "@"

$allData = New-Object 'object[,]' 36,9
$allData[0,0] = 0
$allData[0,1] = '开始时间'
$allData[0,2] = '名称'
$allData[0,3] = '地点'
$allData[0,4] = '具体时间范围'
$allData[0,5] = '想去人数'
$allData[0,6] = '最低票价'
$allData[0,7] = 'Link'
$allData[0,8] = 'Cover'
$allData[1,0] = 1
$allData[1,1] = '2024-06-15'
$allData[1,2] = '上饶·宅舞联萌·随舞动漫派对（免费活动)'
$allData[1,3] = '春江北大道和吉阳路交汇处 槠溪时光PARK'
$allData[1,4] = '2024.06.15 08:00-06.15 21:00'
$allData[1,5] = 36
$allData[1,6] = 22.33
$allData[1,7] = 'https://show.bilibili.com/platform/detail.html?id=85607'
$allData[1,8] = '//i0.hdslb.com/bfs/openplatform/202405/jcZGKqhx1715589649770.jpeg'
$allData[2,0] = 2
$allData[2,1] = '2024-06-22'
$allData[2,2] = '景德镇·BM次元盛典运动番only'
$allData[2,3] = '广场南路金幕影城旁 罗曼园宴会酒店'
$allData[2,4] = '2024.06.22 10:00-06.22 17:00'
$allData[2,5] = 197
$allData[2,6] = 55
$allData[2,7] = 'https://show.bilibili.com/platform/detail.html?id=85197'
$allData[2,8] = '//i2.hdslb.com/bfs/openplatform/202404/Z6eXz0su1714292081978.png'
$allData[3,0] = 3
$allData[3,1] = '2024-06-22'
$allData[3,2] = '萍乡·AU9夏至国漫展'
$allData[3,3] = '金陵东路18号 萍乡市体育馆'
$allData[3,4] = '2024.06.22 10:00-06.22 17:00'
$allData[3,5] = 50
$allData[3,6] = 45
$allData[3,7] = 'https://show.bilibili.com/platform/detail.html?id=86453'
$allData[3,8] = '//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg'
$allData[4,0] = 4
$allData[4,1] = '2024-06-23'
$allData[4,2] = '上饶·BM次元盛典运动番only'
$allData[4,3] = '春江北大道时光PARK内 博悦宴会艺术中心'
$allData[4,4] = '2024.06.23 10:00-06.23 17:00'
$allData[4,5] = 250
$allData[4,6] = 55
$allData[4,7] = 'https://show.bilibili.com/platform/detail.html?id=85201'
$allData[4,8] = '//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png'
$allData[5,0] = 5
$allData[5,1] = '2024-06-23'
$allData[5,2] = '赣州·清风霁月·光夜only'
$allData[5,3] = '平安大道 麋鹿LiveHouse'
$allData[5,4] = '2024.06.23 14:00-06.23 20:00'
$allData[5,5] = 40
$allData[5,6] = 158
$allData[5,7] = 'https://show.bilibili.com/platform/detail.html?id=86993'
$allData[5,8] = '//i1.hdslb.com/bfs/openplatform/202406/PklWR8EP1717429316070.jpeg'
$allData[6,0] = 6
$allData[6,1] = '2024-06-29'
$allData[6,2] = '南昌·第五人格only'
$allData[6,3] = '高处见美好生活公园 百家喜宴高新店'
$allData[6,4] = '2024.06.29 10:00-06.29 17:00'
$allData[6,5] = 138
$allData[6,6] = 65
$allData[6,7] = 'https://show.bilibili.com/platform/detail.html?id=87043'
$allData[6,8] = '//i0.hdslb.com/bfs/openplatform/202405/zir2PYz81717071721569.jpeg'
$allData[7,0] = 7
$allData[7,1] = '2024-06-29'
$allData[7,2] = '萍乡·BM次元盛典运动番only'
$allData[7,3] = '康庄路3号 萍乡梅园国际大酒店'
$allData[7,4] = '2024.06.29 10:00-06.29 17:00'
$allData[7,5] = 253
$allData[7,6] = 55
$allData[7,7] = 'https://show.bilibili.com/platform/detail.html?id=85192'
$allData[7,8] = '//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png'
$allData[8,0] = 8
$allData[8,1] = '2024-06-30'
$allData[8,2] = '南昌·ChinastyleCOSPLAY  '
$allData[8,3] = '真君路999号 南昌玛雅乐园'
$allData[8,4] = '2024.06.30 09:30-07.02 17:30'
$allData[8,5] = 12
$allData[8,6] = 65
$allData[8,7] = 'https://show.bilibili.com/platform/detail.html?id=87045'
$allData[8,8] = '//i1.hdslb.com/bfs/openplatform/202405/wajWy7ID1717149642528.jpeg'
$allData[9,0] = 9
$allData[9,1] = '2024-06-30'
$allData[9,2] = '宜春·BM次元盛典运动番only'
$allData[9,3] = '鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$allData[9,4] = '2024.06.30 10:00-06.30 17:00'
$allData[9,5] = 249
$allData[9,6] = 55
$allData[9,7] = 'https://show.bilibili.com/platform/detail.html?id=84636'
$allData[9,8] = '//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'
$allData[10,0] = 10
$allData[10,1] = '2024-07-06'
$allData[10,2] = '南昌·次元星球动漫游戏展'
$allData[10,3] = '龙蟠街666号融创茂1层 融创茂'
$allData[10,4] = '2024.07.06 10:00-07.06 17:00'
$allData[10,5] = 14
$allData[10,6] = '不可售'
$allData[10,7] = 'https://show.bilibili.com/platform/detail.html?id=86405'
$allData[10,8] = '//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg'
$allData[11,0] = 11
$allData[11,1] = '2024-07-06'
$allData[11,2] = '鹰潭·BM次元盛典运动番only'
$allData[11,3] = '体育馆东路2号九小隔壁 忆江南•宴会楼'
$allData[11,4] = '2024.07.06 10:00-07.06 17:00'
$allData[11,5] = 37
$allData[11,6] = 55
$allData[11,7] = 'https://show.bilibili.com/platform/detail.html?id=85997'
$allData[11,8] = '//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'
$allData[12,0] = 12
$allData[12,1] = '2024-07-07'
$allData[12,2] = '赣州·BM次元盛典运动番only'
$allData[12,3] = '米瑞金路2口0号上客天下1楼 上客天下.老虔州'
$allData[12,4] = '2024.07.07 10:00-07.07 17:00'
$allData[12,5] = 25
$allData[12,6] = 55
$allData[12,7] = 'https://show.bilibili.com/platform/detail.html?id=86602'
$allData[12,8] = '//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png'
$allData[13,0] = 13
$allData[13,1] = '2024-07-12'
$allData[13,2] = '新余·2024第三届MG动漫嘉年华'
$allData[13,3] = '仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅'
$allData[13,4] = '2024.07.12 10:00-07.13 17:30'
$allData[13,5] = 86
$allData[13,6] = 55
$allData[13,7] = 'https://show.bilibili.com/platform/detail.html?id=86536'
$allData[13,8] = '//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg'
$allData[14,0] = 14
$allData[14,1] = '2024-07-13'
$allData[14,2] = '南昌·SuperComic动漫游戏博览会'
$allData[14,3] = '怀玉山大道1315号 南昌绿地国际博览中心'
$allData[14,4] = '2024.07.13 09:00-07.14 17:00'
$allData[14,5] = 433
$allData[14,6] = 65
$allData[14,7] = 'https://show.bilibili.com/platform/detail.html?id=86992'
$allData[14,8] = '//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg'
$allData[15,0] = 15
$allData[15,1] = '2024-07-13'
$allData[15,2] = '宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$allData[15,3] = '宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$allData[15,4] = '2024.07.13 10:00-07.14 17:00'
$allData[15,5] = 44
$allData[15,6] = 55
$allData[15,7] = 'https://show.bilibili.com/platform/detail.html?id=86667'
$allData[15,8] = '//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg'
$allData[16,0] = 16
$allData[16,1] = '2024-07-14'
$allData[16,2] = '吉安·COMIC LIFE次元假日05'
$allData[16,3] = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$allData[16,4] = '2024.07.14 09:00-07.14 18:00'
$allData[16,5] = 476
$allData[16,6] = 52.1
$allData[16,7] = 'https://show.bilibili.com/platform/detail.html?id=85924'
$allData[16,8] = '//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'
$allData[17,0] = 17
$allData[17,1] = '2024-07-19'
$allData[17,2] = '赣州·第四届赣州半夏动漫展'
$allData[17,3] = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$allData[17,4] = '2024.07.19 10:00-07.21 17:00'
$allData[17,5] = 402
$allData[17,6] = 55
$allData[17,7] = 'https://show.bilibili.com/platform/detail.html?id=86587'
$allData[17,8] = '//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg'
$allData[18,0] = 18
$allData[18,1] = '2024-07-20'
$allData[18,2] = '南昌·漫拥动漫嘉年华Pro-追光启航'
$allData[18,3] = '小蓝南路420号 洪州体育馆'
$allData[18,4] = '2024.07.20 09:00-07.21 17:00'
$allData[18,5] = 137
$allData[18,6] = 52.5
$allData[18,7] = 'https://show.bilibili.com/platform/detail.html?id=85796'
$allData[18,8] = '//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'
$allData[19,0] = 19
$allData[19,1] = '2024-07-21'
$allData[19,2] = '乐平·CY境界次元动漫夏时庆'
$allData[19,3] = '翥山西路182号 佳佳基大酒店'
$allData[19,4] = '2024.07.21 10:00-07.21 17:00'
$allData[19,5] = 64
$allData[19,6] = 30
$allData[19,7] = 'https://show.bilibili.com/platform/detail.html?id=86768'
$allData[19,8] = '//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png'
$allData[20,0] = 20
$allData[20,1] = '2024-07-21'
$allData[20,2] = '九江·SXD动漫嘉年华'
$allData[20,3] = '湓浦街道大中路339号 百嘉洲际酒店'
$allData[20,4] = '2024.07.21 10:00-07.21 17:30'
$allData[20,5] = 32
$allData[20,6] = 45
$allData[20,7] = 'https://show.bilibili.com/platform/detail.html?id=86832'
$allData[20,8] = '//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg'
$allData[21,0] = 21
$allData[21,1] = '2024-07-21'
$allData[21,2] = '萍乡·NL14动漫游戏展·夏日狂想曲'
$allData[21,3] = '公园南路168号(近工行城北分理处) 梅生嘉华酒店'
$allData[21,4] = '2024.07.21 10:00-07.21 17:00'
$allData[21,5] = 38
$allData[21,6] = 40
$allData[21,7] = 'https://show.bilibili.com/platform/detail.html?id=86658'
$allData[21,8] = '//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg'
$allData[22,0] = 22
$allData[22,1] = '2024-07-26'
$allData[22,2] = '南昌·萌卡动漫展'
$allData[22,3] = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$allData[22,4] = '2024.07.26 09:00-07.28 17:00'
$allData[22,5] = 1099
$allData[22,6] = 19.9
$allData[22,7] = 'https://show.bilibili.com/platform/detail.html?id=86776'
$allData[22,8] = '//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg'
$allData[23,0] = 23
$allData[23,1] = '2024-07-27'
$allData[23,2] = '江西·次元星河动漫游戏嘉年华'
$allData[23,3] = '九龙大道1177号 南昌绿地国际博览中心'
$allData[23,4] = '2024.07.27 10:00-07.28 17:00'
$allData[23,5] = 2832
$allData[23,6] = 69
$allData[23,7] = 'https://show.bilibili.com/platform/detail.html?id=85493'
$allData[23,8] = '//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png'
$allData[24,0] = 24
$allData[24,1] = '2024-07-27'
$allData[24,2] = '赣州·马娘only'
$allData[24,3] = '火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)'
$allData[24,4] = '2024.07.27 09:00-07.27 17:00'
$allData[24,5] = 22
$allData[24,6] = 60
$allData[24,7] = 'https://show.bilibili.com/platform/detail.html?id=86772'
$allData[24,8] = '//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png'
$allData[25,0] = 25
$allData[25,1] = '2024-07-28'
$allData[25,2] = '赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$allData[25,3] = '兴国路恒大帝景西门 江西长庚控股有限公司'
$allData[25,4] = '2024.07.28 11:00-07.28 17:00'
$allData[25,5] = 55
$allData[25,6] = 56
$allData[25,7] = 'https://show.bilibili.com/platform/detail.html?id=85688'
$allData[25,8] = '//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'
$allData[26,0] = 26
$allData[26,1] = '2024-08-03'
$allData[26,2] = '南昌·幻梦境国际动漫游戏嘉年华1th'
$allData[26,3] = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$allData[26,4] = '2024.08.03 09:00-08.04 17:30'
$allData[26,5] = 536
$allData[26,6] = 64
$allData[26,7] = 'https://show.bilibili.com/platform/detail.html?id=83980'
$allData[26,8] = '//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'
$allData[27,0] = 27
$allData[27,1] = '2024-08-03'
$allData[27,2] = '吉安·COMIC LIFE周年庆典'
$allData[27,3] = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$allData[27,4] = '2024.08.03 09:30-08.03 18:00'
$allData[27,5] = 46
$allData[27,6] = 9.9
$allData[27,7] = 'https://show.bilibili.com/platform/detail.html?id=87164'
$allData[27,8] = '//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg'
$allData[28,0] = 28
$allData[28,1] = '2024-08-03'
$allData[28,2] = '景德镇·第十五届瓷都ACG动漫游戏博览会'
$allData[28,3] = '迎宾大道与寺山路交叉口东200米 陶博城'
$allData[28,4] = '2024.08.03 09:00-08.04 17:00'
$allData[28,5] = 1608
$allData[28,6] = 55
$allData[28,7] = 'https://show.bilibili.com/platform/detail.html?id=86341'
$allData[28,8] = '//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'
$allData[29,0] = 29
$allData[29,1] = '2024-08-03'
$allData[29,2] = '景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$allData[29,3] = '迎宾大道与寺山路交叉口东200米 陶博城'
$allData[29,4] = '2024.08.03 08:30-08.03 17:00'
$allData[29,5] = 566
$allData[29,6] = '已售罄'
$allData[29,7] = 'https://show.bilibili.com/platform/detail.html?id=85981'
$allData[29,8] = '//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'
$allData[30,0] = 30
$allData[30,1] = '2024-08-03'
$allData[30,2] = '樟树·第二届静卿国风动漫文化展览会'
$allData[30,3] = '杏佛路89号 樟树银河国际酒店'
$allData[30,4] = '2024.08.03 09:00-08.03 17:00'
$allData[30,5] = 452
$allData[30,6] = 45
$allData[30,7] = 'https://show.bilibili.com/platform/detail.html?id=86683'
$allData[30,8] = '//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg'
$allData[31,0] = 31
$allData[31,1] = '2024-08-04'
$allData[31,2] = '九江·第一届异次元动漫嘉年华'
$allData[31,3] = '长虹西大道兴城广场99号 九江半岛宾馆'
$allData[31,4] = '2024.08.04 08:00-08.04 17:00'
$allData[31,5] = 261
$allData[31,6] = 45
$allData[31,7] = 'https://show.bilibili.com/platform/detail.html?id=84407'
$allData[31,8] = '//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg'
$allData[32,0] = 32
$allData[32,1] = '2024-08-06'
$allData[32,2] = '南昌·第一届异次元动漫嘉年华'
$allData[32,3] = '民德路411号 东方豪景花园酒店(民德路店)'
$allData[32,4] = '2024.08.06 08:00-08.06 17:00'
$allData[32,5] = 390
$allData[32,6] = 55
$allData[32,7] = 'https://show.bilibili.com/platform/detail.html?id=84102'
$allData[32,8] = '//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'
$allData[33,0] = 33
$allData[33,1] = '2024-08-06'
$allData[33,2] = '宜春·第三十五届静卿国风动漫文化展览会'
$allData[33,3] = '宜阳大道19号(交通银行旁) 宜春安缦文华酒店'
$allData[33,4] = '2024.08.06 09:00-08.06 17:00'
$allData[33,5] = 450
$allData[33,6] = 45
$allData[33,7] = 'https://show.bilibili.com/platform/detail.html?id=86684'
$allData[33,8] = '//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg'
$allData[34,0] = 34
$allData[34,1] = '2024-08-08'
$allData[34,2] = '赣州·第二届异次元动漫嘉年华'
$allData[34,3] = '金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$allData[34,4] = '2024.08.08 08:00-08.08 17:00'
$allData[34,5] = 600
$allData[34,6] = 45
$allData[34,7] = 'https://show.bilibili.com/platform/detail.html?id=84184'
$allData[34,8] = '//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'
$allData[35,0] = 35
$allData[35,1] = '2024-08-10'
$allData[35,2] = '高安·第二届静卿国风动漫文化展览会'
$allData[35,3] = '华林中路606号 高安华鼎国际大酒店'
$allData[35,4] = '2024.08.10 09:00-08.10 17:00'
$allData[35,5] = 422
$allData[35,6] = 45
$allData[35,7] = 'https://show.bilibili.com/platform/detail.html?id=86682'
$allData[35,8] = '//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg'
$wsAll.Range("A1:I36").Value = $allData

$wsAll.Range("B1:B36").ClearFormats()

